{"js": "// Add a new bulleted list item \"Emitir cupom fiscal ap\u00f3s pagamento\" right\n// after the existing \"Consulta de ve\u00edculos estacionados\" list item.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text === \"Consulta de ve\u00edculos estacionados\"\n);\n\nif (!target) {\n  throw new Error('Could not find paragraph \"Consulta de ve\u00edculos estacionados\".');\n}\n\n// Inserting directly after the target paragraph copies its paragraph\n// formatting (ListParagraph style + bullet numbering), matching the target\n// list item's look.\ntarget.insertParagraph(\"Emitir cupom fiscal ap\u00f3s pagamento\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the existing bullet \"Consulta de ve\u00edculos estacionados\".\n$found = $d.Content\n$found.Find.ClearFormatting()\n$found.Find.Execute(\"Consulta de ve\u00edculos estacionados\") | Out-Null\n\nif (-not $found.Find.Found) {\n    throw 'Could not find paragraph \"Consulta de ve\u00edculos estacionados\".'\n}\n\n$target = $found.Paragraphs(1)\n\n# Insert a new paragraph right after it; it inherits the ListParagraph /\n# bullet-numbering formatting from the paragraph it follows.\n$target.Range.InsertParagraphAfter()\n\n$newPara = $target.Next()\n$newPara.Range.Text = \"Emitir cupom fiscal ap\u00f3s pagamento\"\n"}
